$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new pump initial-setup test rows (3-6) that were previously
# blank placeholders (only carried formatting, no values).
$data = @(
    @(2, 1, "F.0.1.22_2", 1440, 5.6, 6.8, 10, 1),
    @(3, 1, "F.0.1.22_3", 1440, 5.6, 6.8, 15, 1),
    @(4, 1, "F.0.1.22_4", 1440, 5.6, 6.8, 3, 1),
    @(5, 1, "F.0.1.21_1", 1440, 5.6, 6.8, 5, 1)
)

$row = 3
foreach ($r in $data) {
    $ws.Range("A$row").Value = $r[0]
    $ws.Range("B$row").Value = $r[1]
    $ws.Range("C$row").Value = $r[2]
    $ws.Range("D$row").Value = $r[3]
    $ws.Range("E$row").Value = $r[4]
    $ws.Range("F$row").Value = $r[5]
    $ws.Range("G$row").Value = $r[6]
    $ws.Range("H$row").Value = $r[7]
    $row++
}

# Move the active selection to H6, matching where editing finished.
[void]$ws.Range("H6").Select()
